$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.215.38"
$ws.Range("E2").Value = "  +0.06%  "
$ws.Range("D3").Value = "1.871.50"
$ws.Range("E3").Value = "  +1.05%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.60"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.35%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4696"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2847"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.47%  "
$ws.Range("B9").Value = "OKB"
$ws.Range("C9").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "41.50"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.48%  "
$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06563"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.11%  "
$ws.Range("B11").Value = "Solana"
$ws.Range("C11").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.23"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.85%  "
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07786"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.23%  "
$ws.Range("B13").Value = "Litecoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "96.29"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.29%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.873.51"
$ws.Range("E14").Value = "  +1.18%  "
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6899"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.16%  "
$ws.Range("B16").Value = "Polkadot"
$ws.Range("C16").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.087"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.16%  "
$ws.Range("B17").Value = "BitcoinCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "267.00"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.61%  "
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "30.199.90"
$ws.Range("E18").Value = "  +0.16%  "
$ws.Range("B19").Value = "Avalanche"
$ws.Range("C19").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.69"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.46%  "
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007713"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.76%  "
$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.001"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.04%  "
$ws.Range("B22").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C22").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D22").Value = "2.114.62"
$ws.Range("E22").Value = "  +0.97%  "
$ws.Range("B23").Value = "BinanceUSD"
$ws.Range("C23").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9997"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("B24").Value = "Uniswap"
$ws.Range("C24").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.234"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.60%  "
$ws.Range("B25").Value = "Chainlink"
$ws.Range("C25").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "6.162"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.37%  "
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.487"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.61%  "
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "165.39"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.19%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.73"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.46%  "
$ws.Range("B29").Value = "LidoDAOToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.936"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.31%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.372"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.57%  "
$ws.Range("B31").Value = "Stellar"
$ws.Range("C31").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09905"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.63%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.358"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.72%  "
$ws.Range("B33").Value = "PancakeSwap"
$ws.Range("C33").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.457"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.50%  "
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.048"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.38%  "
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.04730"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.76%  "
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.128"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.07%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6995"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.17%  "
$ws.Range("B38").Value = "HuobiToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.717"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.43%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01863"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.25%  "
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.779"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.79%  "
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.240"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.31%  "
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "72.71"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.52%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.937"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.24%  "
$ws.Range("B44").Value = "PaxDollar"
$ws.Range("C44").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.000"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.05%  "
$ws.Range("B45").Value = "TheSandbox"
$ws.Range("C45").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4147"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.41%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "102.93"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.13%  "
$ws.Range("B47").Value = "TrustWalletToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.8317"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.76%  "
$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "973.44"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.81%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.063"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.63%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.129"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.03%  "
$ws.Range("B51").Value = "Elrond"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "34.44"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.78%  "
